$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'235.82"

# Row 3
$ws.Cells.Item(3, 4).Value = "'22.34"

# Row 4
$ws.Cells.Item(4, 4).Value = "'5.415"

# Row 5
$ws.Cells.Item(5, 4).Value = "'0.05633"

# Row 6
$ws.Cells.Item(6, 4).Value = "'3.371"

# Row 7
$ws.Cells.Item(7, 4).Value = "'6.478"

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.7833"

# Row 10
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Value = "'0.1396"
$ws.Cells.Item(10, 5).Value = "9WazirXWRX"

# Row 11
$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).Value = "'0.07431"
$ws.Cells.Item(11, 5).Value = "10MandalaExchangeTokenMDX"

# Row 12
$ws.Cells.Item(12, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(12, 4).Value = "'0.03194"
$ws.Cells.Item(12, 5).Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).Value = "'0.02941"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

# Row 14
$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 4).Value = "'0.09262"
$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"

# Row 15
$ws.Cells.Item(15, 2).Value = "BitForexToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15, 4).Value = "'0.001676"
$ws.Cells.Item(15, 5).Value = "14BitForexTokenBF"

# Row 16
$ws.Cells.Item(16, 2).Value = "MCDex"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(16, 4).Value = "'3.260"
$ws.Cells.Item(16, 5).Value = "15MCDexMCB"

# Row 17
$ws.Cells.Item(17, 2).Value = "CoinExToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(17, 4).Value = "'0.04753"
$ws.Cells.Item(17, 5).Value = "16CoinExTokenCET"

# Row 18
$ws.Cells.Item(18, 2).Value = "One"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(18, 4).Value = "'0.0005740"
$ws.Cells.Item(18, 5).Value = "17OneONE"

# Row 20
$ws.Cells.Item(20, 4).Value = "'0.005112"

# Row 21
$ws.Cells.Item(21, 4).Value = "'0.001050"

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.0001500"

# Row 23
$ws.Cells.Item(23, 4).Value = "'3.894"

# Row 24
$ws.Cells.Item(24, 4).Value = "'2.146"

# Row 27
$ws.Cells.Item(27, 4).Value = "'0.0004990"

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.04051"

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.007005"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"

# Row 42
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42, 4).Value = "'0.1041"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"

# Row 43
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(43, 4).Value = "'0.003371"
$ws.Cells.Item(43, 5).Value = "42CEJICEJIBestin24h"

# Row 44
$ws.Cells.Item(44, 4).Value = "'0.009300"

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.00005435"

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.6752"

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.03956"

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.00002100"
